# This script applies updated TP/FP/FN/Total/Precision/Recall figures to the
# class_summary_df worksheet, reflecting the vote-validation module results
# obtained after working on the balanced dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Row, C=TP, D=FP, E=FN, F=Total, G=Precision, H=Recall
$data = @(
    @(2, 0, 47, 99, 146, 0, 0),
    @(3, 84, 15, 14, 113, 0.8484848484848485, 0.8571428571428571),
    @(4, 71, 28, 27, 126, 0.7171717171717171, 0.7244897959183674),
    @(5, 95, 4, 4, 103, 0.9595959595959596, 0.9595959595959596),
    @(6, 88, 11, 11, 110, 0.8888888888888888, 0.8888888888888888),
    @(7, 86, 13, 13, 112, 0.8686868686868687, 0.8686868686868687),
    @(8, 74, 25, 20, 119, 0.7474747474747475, 0.7872340425531915),
    @(9, 80, 19, 17, 116, 0.8080808080808081, 0.8247422680412371),
    @(10, 94, 5, 4, 103, 0.9494949494949495, 0.9591836734693877),
    @(11, 99, 0, 0, 99, 1, 1),
    @(12, 80, 19, 17, 116, 0.8080808080808081, 0.8247422680412371),
    @(13, 93, 6, 6, 105, 0.9393939393939394, 0.9393939393939394),
    @(14, 65, 34, 28, 127, 0.6565656565656566, 0.6989247311827957),
    @(15, 94, 5, 5, 104, 0.9494949494949495, 0.9494949494949495),
    @(16, 48, 3, 3, 54, 0.9411764705882353, 0.9411764705882353),
    @(17, 84, 15, 12, 111, 0.8484848484848485, 0.875),
    @(18, 83, 16, 14, 113, 0.8383838383838383, 0.8556701030927835),
    @(19, 84, 15, 13, 112, 0.8484848484848485, 0.865979381443299),
    @(20, 86, 13, 10, 109, 0.8686868686868687, 0.8958333333333334),
    @(21, 80, 19, 19, 118, 0.8080808080808081, 0.8080808080808081),
    @(22, 91, 8, 7, 106, 0.9191919191919192, 0.9285714285714286),
    @(23, 95, 4, 4, 103, 0.9595959595959596, 0.9595959595959596),
    @(24, 66, 33, 25, 124, 0.6666666666666666, 0.7252747252747253),
    @(25, 77, 22, 16, 115, 0.7777777777777778, 0.8279569892473119),
    @(26, 81, 18, 17, 116, 0.8181818181818182, 0.826530612244898),
    @(27, 81, 18, 12, 111, 0.8181818181818182, 0.8709677419354839),
    @(28, 72, 27, 24, 123, 0.7272727272727273, 0.75),
    @(29, 56, 43, 30, 129, 0.5656565656565656, 0.6511627906976745),
    @(30, 78, 21, 13, 112, 0.7878787878787878, 0.8571428571428571),
    @(31, 88, 11, 9, 108, 0.8888888888888888, 0.9072164948453608),
    @(32, 89, 10, 7, 106, 0.898989898989899, 0.9270833333333334),
    @(33, 96, 3, 3, 102, 0.9696969696969697, 0.9696969696969697),
    @(34, 72, 27, 23, 122, 0.7272727272727273, 0.7578947368421053),
    @(35, 93, 6, 5, 104, 0.9393939393939394, 0.9489795918367347),
    @(36, 82, 17, 14, 113, 0.8282828282828283, 0.8541666666666666),
    @(37, 82, 17, 11, 110, 0.8282828282828283, 0.8817204301075269),
    @(38, 82, 17, 17, 116, 0.8282828282828283, 0.8282828282828283),
    @(39, 83, 16, 11, 110, 0.8383838383838383, 0.8829787234042553),
    @(40, 84, 15, 15, 114, 0.8484848484848485, 0.8484848484848485),
    @(41, 46, 2, 2, 50, 0.9583333333333334, 0.9583333333333334),
    @(42, 94, 5, 5, 104, 0.9494949494949495, 0.9494949494949495),
    @(43, 80, 19, 17, 116, 0.8080808080808081, 0.8247422680412371),
    @(44, 99, 0, 0, 99, 1, 1),
    @(45, 94, 5, 5, 104, 0.9494949494949495, 0.9494949494949495)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
}
